# Update gh-pages to output generated at 456a3b4
# Applies updated "want-to-go" counts (column F) and minimum ticket price (column G)
# to the "展览" (sheet 1) and "全部类型" (sheet 4) worksheets.

$wb = $excel.ActiveWorkbook

$targetSheets = @($wb.Worksheets.Item(1), $wb.Worksheets.Item(4))

foreach ($ws in $targetSheets) {
    $ws.Range("G2").Value = 55

    $ws.Range("F4").Value = 284
    $ws.Range("F5").Value = 55
    $ws.Range("F8").Value = 2055
    $ws.Range("F9").Value = 71
    $ws.Range("F11").Value = 4477
    $ws.Range("F15").Value = 12
    $ws.Range("F16").Value = 126
    $ws.Range("F18").Value = 19
    $ws.Range("F19").Value = 82
    $ws.Range("F20").Value = 3331
    $ws.Range("F22").Value = 503
    $ws.Range("F25").Value = 82
    $ws.Range("F29").Value = 60
    $ws.Range("F30").Value = 204
    $ws.Range("F32").Value = 634
    $ws.Range("F33").Value = 1975
    $ws.Range("F34").Value = 353
}
